$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update USERID (col G) and PASSWORD (col H) for rows 2-5, plus the
# PREPARATION text in col F which embeds the username/password.

$ws.Range("G2").Value = 32382
$ws.Range("H2").Value = "bni1234"
$ws.Range("F2").Value = "Username : 32382;`nPassword : bni1234;`nKode Level Approval : 9;`nNama Level Approval : Team Lead"

$ws.Range("G3").Value = 32382
$ws.Range("H3").Value = "bni1234"
$ws.Range("F3").Value = "Username : 32382;`nPassword : bni1234;`nKode Level Approval : 9"

$ws.Range("G4").Value = 32382
$ws.Range("H4").Value = "bni1234"
$ws.Range("F4").Value = "Username : 32382;`nPassword : bni1234;`nKode Level Approval : 9;`nNama Level Approval : Team Leader"

$ws.Range("G5").Value = 32382
$ws.Range("H5").Value = "bni1234"
$ws.Range("F5").Value = "Username : 32382;`nPassword : bni1234;`nKode Level Approval : 9"

# Update active cell selection to F2
$ws.Range("F2").Select()
